# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Coin/Link/Price/Volume cells are stored as plain text in the source sheet. Excel's
# COM layer auto-coerces numeric-looking strings (e.g. "1.00", "0.260") into actual
# numbers on assignment, so force a text number format first on any cell whose new
# value would otherwise be reinterpreted, keeping it identical to the scraped text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.929.15'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '1.630.19'
$ws.Range("E3").Value = '  +1.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.56'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.523'
$ws.Range("E6").Value = '  +1.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.32'
$ws.Range("E8").Value = '  +9.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.260'
$ws.Range("E9").Value = '  +3.37%  '
$ws.Range("E10").Value = '  +2.55%  '
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '1.863.90'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("D13").Value = '1.646.49'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.572'
$ws.Range("E14").Value = '  +6.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.91'
$ws.Range("E15").Value = '  +5.25%  '
$ws.Range("D16").Value = '29.952.71'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.98'
$ws.Range("E17").Value = '  +18.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.71'
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.69'
$ws.Range("E19").Value = '  +3.22%  '
$ws.Range("E20").Value = '  +2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  +3.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.61'
$ws.Range("E23").Value = '  +4.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.29'
$ws.Range("E25").Value = '  +2.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.73'
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("E27").Value = '  +2.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.63'
$ws.Range("E28").Value = '  +3.47%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0492'
$ws.Range("E30").Value = '  +3.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.13'
$ws.Range("E31").Value = '  +5.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +4.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.21'
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("D34").Value = '1.430.34'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("E35").Value = '  +6.63%  '
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.90'
$ws.Range("E37").Value = '  +2.46%  '
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.560'
$ws.Range("E40").Value = '  +4.13%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0502'
$ws.Range("E41").Value = '  +2.25%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.834'
$ws.Range("E42").Value = '  +4.33%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.97'
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.60'
$ws.Range("E44").Value = '  +3.16%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.04'
$ws.Range("E45").Value = '  +4.58%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '69.88'
$ws.Range("E46").Value = '  +6.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("D49").Value = '1.772.11'
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '89.45'
$ws.Range("E50").Value = '  +3.55%  '
$ws.Range("D51").Value = '0.0₆0107'
$ws.Range("E51").Value = '  +1.76%  '
